$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force pure-numeric-looking target strings to remain text by applying a text number format
# before assignment (Excel reuses a single shared style for all of these).
$textForceRefs = @("D5", "D8", "D11", "D16", "D18", "D20", "D23", "D25", "D26", "D29", "D31", "D36", "D37", "D40", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D51")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "27.480.34"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "1.618.07"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "211.32"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  -1.48%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "22.85"
$ws.Range("E8").Value = "  -1.59%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").Value = "0.0887"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.846.88"
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").Value = "1.619.74"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("E14").Value = "  -0.31%  "
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").Value = "65.00"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "27.457.08"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "232.73"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "10.17"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +5.50%  "
$ws.Range("D25").Value = "150.54"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("D26").Value = "6.86"
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "0.0483"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").Value = "1.474.73"
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("E34").Value = "  -3.13%  "
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "0.965"
$ws.Range("E36").Value = "  +8.50%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -0.35%  "
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "0.862"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "68.00"
$ws.Range("E42").Value = "  +2.35%  "
$ws.Range("B43").Value = "mCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D43").Value = "2.46"
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "0.982"
$ws.Range("E44").Value = "  -5.07%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "2.20"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "5.27"
$ws.Range("E46").Value = "  -7.04%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.757.46"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "1.74"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "86.81"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0105"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.101"
$ws.Range("E51").Value = "  +1.40%  "
